$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.840.46'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.512.30'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.03'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.47'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.66%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").Value = '2.510.32'
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("D10").ClearFormats()
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = '2.955.05'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.00'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.44%  '
$ws.Range("D16").Value = '58.779.19'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000140'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '2.508.99'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.27'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.04'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.81'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.07'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.424'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.13%  '
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.76'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.68'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").Value = '0.0₃0768'
$ws.Range("E30").Value = '  -1.13%  '
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.46'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.22%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -6.14%  '
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.43'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.23%  '
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.81'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.801'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.20'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '278.42'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.61'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.38'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0510'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("E51").Value = '  -2.57%  '
